$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.668968642179323
$ws.Range("D2").Value = 0.009111524370560176
$ws.Range("E2").Value = 0.06081283538677518
$ws.Range("F2").Value = 4.17909545117071
$ws.Range("G2").Value = 0.002634859792585293
$ws.Range("J2").Value = 0.1478672292744827
$ws.Range("K2").Value = 1.868411174977524
$ws.Range("L2").Value = 0.4569260639090515
$ws.Range("N2").Value = 3.975806658014051
$ws.Range("B3").Value = 1.64734674111898
$ws.Range("D3").Value = 0.008667039786935504
$ws.Range("E3").Value = 0.06057112385611951
$ws.Range("F3").Value = 4.16234374125051
$ws.Range("G3").Value = 0.002639890346336095
$ws.Range("J3").Value = 0.1478545024921321
$ws.Range("K3").Value = 1.802384033067739
$ws.Range("L3").Value = 0.4474931867815144
$ws.Range("N3").Value = 3.985893407974686
$ws.Range("B4").Value = 1.634912823024422
$ws.Range("D4").Value = 0.008391523925883604
$ws.Range("E4").Value = 0.0604205277641805
$ws.Range("F4").Value = 4.153718281887734
$ws.Range("G4").Value = 0.002643143361588096
$ws.Range("J4").Value = 0.1478494338654501
$ws.Range("K4").Value = 1.763006237460445
$ws.Range("L4").Value = 0.4419676822564327
$ws.Range("N4").Value = 3.992924812539641
$ws.Range("B5").Value = 1.63005783711202
$ws.Range("D5").Value = 0.008278564649433662
$ws.Range("E5").Value = 0.06035859748312378
$ws.Range("F5").Value = 4.150620493967494
$ws.Range("G5").Value = 0.00264451042641876
$ws.Range("J5").Value = 0.147848053895455
$ws.Range("K5").Value = 1.747251355115736
$ws.Range("L5").Value = 0.4397829055083378
$ws.Range("N5").Value = 3.996000726232054
$ws.Range("B6").Value = 1.62926447705172
$ws.Range("D6").Value = 0.008259765633120253
$ws.Range("E6").Value = 0.06034827976961887
$ws.Range("F6").Value = 4.150131295292582
$ws.Range("G6").Value = 0.002644739933247458
$ws.Range("J6").Value = 0.147847866016285
$ws.Range("K6").Value = 1.74465288141576
$ws.Range("L6").Value = 0.439424164675188
$ws.Range("N6").Value = 3.996524191205424
$ws.Range("B7").Value = 1.634846488591165
$ws.Range("D7").Value = 0.008390003328845097
$ws.Range("E7").Value = 0.06041969483782283
$ws.Range("F7").Value = 4.15367481540504
$ws.Range("G7").Value = 0.002643161630593893
$ws.Range("J7").Value = 0.1478494124850864
$ws.Range("K7").Value = 1.762792580535802
$ws.Range("L7").Value = 0.4419379467463784
$ws.Range("N7").Value = 3.992965443114727
$ws.Range("B8").Value = 1.661338759933869
$ws.Range("D8").Value = 0.008958787345555663
$ws.Range("E8").Value = 0.06072993933679083
$ws.Range("F8").Value = 4.172974696202587
$ws.Range("G8").Value = 0.00263656032315439
$ws.Range("J8").Value = 0.1478622682402673
$ws.Range("K8").Value = 1.845403221923931
$ws.Range("L8").Value = 0.4536182936195701
$ws.Range("N8").Value = 3.979110506455299
$ws.Range("B9").Value = 1.719967143853438
$ws.Range("D9").Value = 0.01005495431777348
$ws.Range("E9").Value = 0.06132155760864899
$ws.Range("F9").Value = 4.224014630039193
$ws.Range("G9").Value = 0.002624912075569377
$ws.Range("J9").Value = 0.147909493644216
$ws.Range("K9").Value = 2.016669608877123
$ws.Range("L9").Value = 0.4786410938785224
$ws.Range("N9").Value = 3.958600355963057
$ws.Range("B10").Value = 1.76711427995474
$ws.Range("D10").Value = 0.01085062149188332
$ws.Range("E10").Value = 0.0617468062663189
$ws.Range("F10").Value = 4.269594151167723
$ws.Range("G10").Value = 0.002617135999703368
$ws.Range("J10").Value = 0.14795794286213
$ws.Range("K10").Value = 2.14821810852834
$ws.Range("L10").Value = 0.498325471845277
$ws.Range("N10").Value = 3.947605628452948
$ws.Range("B11").Value = 1.789448226814102
$ws.Range("D11").Value = 0.01121093450360888
$ws.Range("E11").Value = 0.06193839777166943
$ws.Range("F11").Value = 4.292093014752396
$ws.Range("G11").Value = 0.002613766380090523
$ws.Range("J11").Value = 0.147983043806895
$ws.Range("K11").Value = 2.20932104202484
$ws.Range("L11").Value = 0.5075649140161573
$ws.Range("N11").Value = 3.943491781894991
$ws.Range("B12").Value = 1.7980329387604
$ws.Range("D12").Value = 0.01134717558364429
$ws.Range("E12").Value = 0.06201069567945261
$ws.Range("F12").Value = 4.300867058578604
$ws.Range("G12").Value = 0.002612514374239814
$ws.Range("J12").Value = 0.1479929950752972
$ws.Range("K12").Value = 2.232641449326763
$ws.Range("L12").Value = 0.5111047494920626
$ws.Range("N12").Value = 3.942061893400066
$ws.Range("B13").Value = 1.796178406563627
$ws.Range("D13").Value = 0.01131784204843456
$ws.Range("E13").Value = 0.06199513609848051
$ws.Range("F13").Value = 4.298966097310654
$ws.Range("G13").Value = 0.00261278295091498
$ws.Range("J13").Value = 0.1479908319568284
$ws.Range("K13").Value = 2.227610873499202
$ws.Range("L13").Value = 0.5103405550222249
$ws.Range("N13").Value = 3.942364150889361
$ws.Range("B14").Value = 1.790151945124677
$ws.Range("D14").Value = 0.01122214698316526
$ws.Range("E14").Value = 0.06194435077049398
$ws.Range("F14").Value = 4.292809763033688
$ws.Range("G14").Value = 0.002613662896701802
$ws.Range("J14").Value = 0.1479838535290359
$ws.Range("K14").Value = 2.211235973045575
$ws.Range("L14").Value = 0.5078553150045479
$ws.Range("N14").Value = 3.943371578433897
$ws.Range("B15").Value = 1.786477141608003
$ws.Range("D15").Value = 0.01116350575509628
$ws.Range("E15").Value = 0.06191321065743338
$ws.Range("F15").Value = 4.289071950928303
$ws.Range("G15").Value = 0.002614205010236627
$ws.Range("J15").Value = 0.1479796373116056
$ws.Range("K15").Value = 2.201229604667731
$ws.Range("L15").Value = 0.5063383835394433
$ws.Range("N15").Value = 3.944005326460953
$ws.Range("B16").Value = 1.765672529766078
$ws.Range("D16").Value = 0.01082704410076829
$ws.Range("E16").Value = 0.0617342490898789
$ws.Range("F16").Value = 4.268159344513407
$ws.Range("G16").Value = 0.002617359577094467
$ws.Range("J16").Value = 0.1479563645558661
$ws.Range("K16").Value = 2.144250322880964
$ws.Range("L16").Value = 0.4977273933339319
$ws.Range("N16").Value = 3.947892365161081
$ws.Range("B17").Value = 1.753136518023098
$ws.Range("D17").Value = 0.01062023991100247
$ws.Range("E17").Value = 0.06162399593047674
$ws.Range("F17").Value = 4.255782406835607
$ws.Range("G17").Value = 0.00261933767747217
$ws.Range("J17").Value = 0.1479428756673529
$ws.Range("K17").Value = 2.109618837870755
$ws.Range("L17").Value = 0.4925178660847678
$ws.Range("N17").Value = 3.950504501718882
$ws.Range("B18").Value = 1.746009583271615
$ws.Range("D18").Value = 0.01050113517302975
$ws.Range("E18").Value = 0.06156040570420362
$ws.Range("F18").Value = 4.24882957432277
$ws.Range("G18").Value = 0.002620491225867835
$ws.Range("J18").Value = 0.1479354049861978
$ws.Range("K18").Value = 2.089818356352168
$ws.Range("L18").Value = 0.4895482963140267
$ws.Range("N18").Value = 3.952090471722741
$ws.Range("B19").Value = 1.743610862325653
$ws.Range("D19").Value = 0.01046078045346377
$ws.Range("E19").Value = 0.06153884459824299
$ws.Range("F19").Value = 4.246503970402784
$ws.Range("G19").Value = 0.002620884514647576
$ws.Range("J19").Value = 0.1479329247891248
$ws.Range("K19").Value = 2.083134606499527
$ws.Range("L19").Value = 0.4885474534635108
$ws.Range("N19").Value = 3.952641793089526
$ws.Range("B20").Value = 1.754462362440762
$ws.Range("D20").Value = 0.01064227051343991
$ws.Range("E20").Value = 0.06163575063702531
$ws.Range("F20").Value = 4.257082763594326
$ws.Range("G20").Value = 0.002619125471098605
$ws.Range("J20").Value = 0.147944281756248
$ws.Range("K20").Value = 2.113293135237654
$ws.Range("L20").Value = 0.4930696535367076
$ws.Range("N20").Value = 3.950217787454974
$ws.Range("B21").Value = 1.791918608311448
$ws.Range("D21").Value = 0.01125026014106822
$ws.Range("E21").Value = 0.06195927442915661
$ws.Range("F21").Value = 4.29461112450582
$ws.Range("G21").Value = 0.002613403785064999
$ws.Range("J21").Value = 0.1479858911082759
$ws.Range("K21").Value = 2.216040731378939
$ws.Range("L21").Value = 0.5085841751835574
$ws.Range("N21").Value = 3.943072198102925
$ws.Range("B22").Value = 1.817140551623027
$ws.Range("D22").Value = 0.01164645743092763
$ws.Range("E22").Value = 0.06216924384559608
$ws.Range("F22").Value = 4.320620054953423
$ws.Range("G22").Value = 0.002609804140023371
$ws.Range("J22").Value = 0.1480156881736523
$ws.Range("K22").Value = 2.284253724634766
$ws.Range("L22").Value = 0.5189631802511769
$ws.Range("N22").Value = 3.939147987005683
$ws.Range("B23").Value = 1.803611267023882
$ws.Range("D23").Value = 0.0114350941133381
$ws.Range("E23").Value = 0.06205730947053034
$ws.Range("F23").Value = 4.306602836504851
$ws.Range("G23").Value = 0.002611712587753641
$ws.Range("J23").Value = 0.147999544770649
$ws.Range("K23").Value = 2.247749788803389
$ws.Range("L23").Value = 0.5134017775998245
$ws.Range("N23").Value = 3.941174071984008
$ws.Range("B24").Value = 1.75386269840061
$ws.Range("D24").Value = 0.01063231113475638
$ws.Range("E24").Value = 0.06163043697181525
$ws.Range("F24").Value = 4.256494365118925
$ws.Range("G24").Value = 0.00261922135864546
$ws.Range("J24").Value = 0.1479436451792635
$ws.Range("K24").Value = 2.11163164400341
$ws.Range("L24").Value = 0.4928201111424784
$ws.Range("N24").Value = 3.950347148648021
$ws.Range("B25").Value = 1.703391570188558
$ws.Range("D25").Value = 0.009760261415376448
$ws.Range("E25").Value = 0.06116322146896769
$ws.Range("F25").Value = 4.208790734649654
$ws.Range("G25").Value = 0.002627925298715834
$ws.Range("J25").Value = 0.1478943271481032
$ws.Range("K25").Value = 1.96933827535247
$ws.Range("L25").Value = 0.4716441206653883
$ws.Range("N25").Value = 3.963434503262761
